# Apply the diff: update the "Fitness" values in column C.
# Rows 2-50 (inclusive) change from 7534/7320 to 7310.
# Rows 51-185 (inclusive) change from 7320/7318/7310 to 7293.
# Rows 186-252 remain unchanged (already 7293).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C50").Value = 7310
$ws.Range("C51:C185").Value = 7293
